$d = $word.ActiveDocument
$wmain = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# =====================================================================
# STEP 1: collapse the three runs "GM " + "-" (en dash) + " " into a
# single run "GM - " in the "GM - Deus ex plasmatio" paragraph.
# =====================================================================
$dash = [char]0x2013
$findRange = $d.Content
[void]$findRange.Find.Execute("GM", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeStart = $findRange.Start
$mergeEnd = $mergeStart + 5
$mergeRange = $d.Range($mergeStart, $mergeEnd)
# two-phase set: first to a same-length sentinel so the run-merge
# actually takes effect, then to the real text
$mergeRange.Text = "ZZZZZ"
$mergeRange2 = $d.Range($mergeStart, $mergeStart + 5)
$mergeRange2.Text = "GM " + $dash + " "

# =====================================================================
# STEP 2: drop the paragraph-mark-only "da-DK" language formatting on
# the "Good - Deus ex bonis" paragraph (its <w:pPr><w:rPr><w:lang.../>).
# =====================================================================
$goodParaIndex = 8
$goodPara = $d.Paragraphs.Item($goodParaIndex)
$goodFull = $d.Range($goodPara.Range.Start, $goodPara.Range.End)
$goodXml = "<w:p $wmain><w:r><w:t>Good</w:t></w:r><w:r><w:t xml:space=`"preserve`"> $dash </w:t></w:r><w:r><w:rPr><w:lang w:val=`"la-Latn`"/></w:rPr><w:t xml:space=`"preserve`">Deus </w:t></w:r><w:r><w:rPr><w:lang w:val=`"la-Latn`"/></w:rPr><w:t>ex bonis</w:t></w:r></w:p><w:p $wmain/>"
$goodFull.InsertXML($goodXml)
$strayPara = $d.Paragraphs.Item($goodParaIndex + 1)
$strayPara.Range.Delete()

# =====================================================================
# STEP 3: rework the end of the "Artefacts:" section.
#   - strip the _GoBack bookmark from its current (empty) paragraph
#     right after "Artefact building"
#   - append, after "Corrupting artefacts":
#       * a blank paragraph
#       * a "Big Bad of first campaign" Heading-1 paragraph
#       * a paragraph with the BBEG idea text, carrying the _GoBack
#         bookmark at its end
# =====================================================================
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfDoc = $lastPara.Range.End - 1
$insertPoint = $d.Range($endOfDoc, $endOfDoc)

$bbegText = "Keeps restarting when he dies, in the start tried to help everyone, in the end he became corrupted by the repeated restarts, leading to logical leaps making him think everyone can be saved if he controls magic and the lay-lines."

$newXml = "<w:p $wmain/>" +
          "<w:p $wmain><w:pPr><w:pStyle w:val=`"Overskrift1`"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Big Bad of first campaign</w:t></w:r></w:p>" +
          "<w:p $wmain><w:r><w:t>$bbegText</w:t></w:r></w:p>"
$insertPoint.InsertXML($newXml)

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPos = $finalPara.Range.End - 1
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
[void]$d.Bookmarks.Add("_GoBack", $bookmarkRange)
